# Generate Report for Handback
# Refresh the timestamps that were recorded when the handback report was generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the a7dc4ec9... source file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-18 17:04:22"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-18 17:04:17"
$wsZhCn.Range("K2").Value = "2016-08-18 17:04:46"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-18 17:04:22"
$wsDeDe.Range("K2").Value = "2016-08-18 17:04:54"
